$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update sheet1 ("irrigation_log") headers ---
$ws1.Range("B1").Value = "Corn IRR (hr)"
$ws1.Range("D1").Value = "Corn FERT (hr)"
$ws1.Range("C1").Value = "Soybean IRR (hr)"
$ws1.Range("E1").Value = "Soybean FERT (hr)"

$ws1.Cells.Item(2, 2).Value = 5
$ws1.Cells.Item(2, 3).Value = 5
$ws1.Cells.Item(2, 4).Value = 5
$ws1.Cells.Item(2, 5).Value = 5
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 3).Value = 5
$ws1.Cells.Item(3, 4).Value = 1
$ws1.Cells.Item(3, 5).Value = 5
$ws1.Cells.Item(4, 2).Value = 5
$ws1.Cells.Item(4, 3).Value = 5
$ws1.Cells.Item(4, 4).Value = 5
$ws1.Cells.Item(4, 5).Value = 5
$ws1.Cells.Item(5, 2).Value = 2
$ws1.Cells.Item(5, 3).Value = 2
$ws1.Cells.Item(5, 4).Value = 2
$ws1.Cells.Item(5, 5).Value = 5
$ws1.Cells.Item(6, 2).Value = 2
$ws1.Cells.Item(6, 3).Value = 2
$ws1.Cells.Item(6, 4).Value = 4
$ws1.Cells.Item(6, 5).Value = 5
$ws1.Cells.Item(7, 2).Value = 5
$ws1.Cells.Item(7, 3).Value = 1
$ws1.Cells.Item(7, 4).Value = 3
$ws1.Cells.Item(7, 5).Value = 1
$ws1.Cells.Item(8, 2).Value = 5
$ws1.Cells.Item(8, 3).Value = 4
$ws1.Cells.Item(8, 4).Value = 3
$ws1.Cells.Item(8, 5).Value = 5
$ws1.Cells.Item(9, 2).Value = 2
$ws1.Cells.Item(9, 3).Value = 5
$ws1.Cells.Item(9, 4).Value = 3
$ws1.Cells.Item(9, 5).Value = 2
$ws1.Cells.Item(10, 2).Value = 1
$ws1.Cells.Item(10, 3).Value = 1
$ws1.Cells.Item(10, 4).Value = 1
$ws1.Cells.Item(10, 5).Value = 5
$ws1.Cells.Item(11, 2).Value = 5
$ws1.Cells.Item(11, 3).Value = 5
$ws1.Cells.Item(11, 4).Value = 5
$ws1.Cells.Item(11, 5).Value = 2
$ws1.Cells.Item(12, 2).Value = 5
$ws1.Cells.Item(12, 3).Value = 1
$ws1.Cells.Item(12, 4).Value = 4
$ws1.Cells.Item(12, 5).Value = 5
$ws1.Cells.Item(13, 2).Value = 2
$ws1.Cells.Item(13, 3).Value = 3
$ws1.Cells.Item(13, 4).Value = 3
$ws1.Cells.Item(13, 5).Value = 2
$ws1.Cells.Item(14, 2).Value = 2
$ws1.Cells.Item(14, 3).Value = 5
$ws1.Cells.Item(14, 4).Value = 1
$ws1.Cells.Item(14, 5).Value = 5
$ws1.Cells.Item(15, 2).Value = 2
$ws1.Cells.Item(15, 3).Value = 3
$ws1.Cells.Item(15, 4).Value = 5
$ws1.Cells.Item(15, 5).Value = 5
$ws1.Cells.Item(16, 2).Value = 5
$ws1.Cells.Item(16, 3).Value = 5
$ws1.Cells.Item(16, 4).Value = 1
$ws1.Cells.Item(16, 5).Value = 4
$ws1.Cells.Item(17, 2).Value = 5
$ws1.Cells.Item(17, 3).Value = 2
$ws1.Cells.Item(17, 4).Value = 5
$ws1.Cells.Item(17, 5).Value = 5
$ws1.Cells.Item(18, 2).Value = 5
$ws1.Cells.Item(18, 3).Value = 5
$ws1.Cells.Item(18, 4).Value = 1
$ws1.Cells.Item(18, 5).Value = 3
$ws1.Cells.Item(19, 2).Value = 1
$ws1.Cells.Item(19, 3).Value = 3
$ws1.Cells.Item(19, 4).Value = 5
$ws1.Cells.Item(19, 5).Value = 2
$ws1.Cells.Item(20, 2).Value = 5
$ws1.Cells.Item(20, 3).Value = 5
$ws1.Cells.Item(20, 4).Value = 5
$ws1.Cells.Item(20, 5).Value = 5
$ws1.Cells.Item(21, 2).Value = 5
$ws1.Cells.Item(21, 3).Value = 5
$ws1.Cells.Item(21, 4).Value = 5
$ws1.Cells.Item(21, 5).Value = 5

# --- Selection state on sheet1 ---
$ws1.Range("B17").Select()

# --- Add the new "deficit_adj" sheet after irrigation_log ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "deficit_adj"

$ws2.Range("A1").Value = "Date"
$ws2.Range("B1").Value = "Corn IRR (inch)"
$ws2.Range("C1").Value = "Soybean IRR (inch)"
$ws2.Range("D1").Value = "Corn FERT (inch)"
$ws2.Range("E1").Value = "Soybean FERT (inch)"

$ws2.Range("A2").Value = [DateTime]::new(2024,4,29)
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 2
$ws2.Range("D2").Value = "-"
$ws2.Range("E2").Value = 1

$ws2.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws2.Columns.Item(4).ColumnWidth = 13.830729166666666
$ws2.Columns.Item(5).ColumnWidth = 16.830729166666668

$ws2.Activate()
$ws2.Range("F2").Select()
